$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Fill in the previously-empty "Te Yang" evaluation block (rows 8-12)
# ---------------------------------------------------------------------------
$teYang = @{
    8  = @(1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,3,3,3,2.5,3,2.5)
    9  = @(1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,3,3,2.5,3,2.5,2.5)
    10 = @(1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,3,2.5,3,2.5,3.5,2.5)
    11 = @(1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,2,2.5,2.5,3.5,3.5,2.5)
    12 = @(1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,3,3,3,3,3.5,2.5)
}

$cols = @("C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z")

foreach ($row in $teYang.Keys) {
    $values = $teYang[$row]
    for ($i = 0; $i -lt $cols.Count; $i++) {
        $ws.Range("$($cols[$i])$row").Value = $values[$i]
    }
}

# ---------------------------------------------------------------------------
# 2) Add a new "Average" block (rows 28-32) averaging each metric across the
#    five judges (rows 3/8/13/18/23, 4/9/14/19/24, ... 7/12/17/22/27)
# ---------------------------------------------------------------------------
$ws.Range("A28").Value = "Average"
$ws.Range("B28").Value = "Fluency"
$ws.Range("B29").Value = "Coherence"
$ws.Range("B30").Value = "Meaningfulness"
$ws.Range("B31").Value = "Poeticness"
$ws.Range("B32").Value = "Overall"

$judgeRowOffsets = @(3, 8, 13, 18, 23)
for ($k = 0; $k -lt 5; $k++) {
    $targetRow = 28 + $k
    foreach ($col in $cols) {
        $refs = @()
        foreach ($base in $judgeRowOffsets) {
            $refs += "$col$($base + $k)"
        }
        $formula = "=AVERAGE(" + ($refs -join ",") + ")"
        $ws.Range("$col$targetRow").Formula = $formula
    }
}

# ---------------------------------------------------------------------------
# 3) Formatting: merge the new label column (A28:A32) and center +
#    vertically center its text, matching the style used for "Average".
# ---------------------------------------------------------------------------
$ws.Range("A28:A32").MergeCells = $true
$ws.Range("A28:A32").HorizontalAlignment = -4108  # xlCenter
$ws.Range("A28:A32").VerticalAlignment = -4108    # xlCenter

# ---------------------------------------------------------------------------
# 4) View tweaks: zoom level and the active selection in the frozen pane
# ---------------------------------------------------------------------------
$win = $excel.ActiveWindow
$win.Zoom = 130
$ws.Range("C35").Select() | Out-Null

$wb.Save()
